# Commit: "git push origin master"
# Adds two new tracked LeetCode questions (rows 47 & 48) to the
# "Questions Tracking Sheet" worksheet, mirroring the formatting of the
# last existing row (46), wires up the new hyperlinks for the "Problem
# link" column, resizes the two new rows, fixes up the used range and
# moves the sheet's scroll/selection to show the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions Tracking Sheet")

# ---------------------------------------------------------------------
# Row 47: "Sum of Nodes with Even-Valued Grandparent"
# ---------------------------------------------------------------------
$ws.Range("A47").Value = 44096
$ws.Range("B47").Value = 1315
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = "https://leetcode.com/problems/sum-of-nodes-with-even-valued-grandparent/"
$ws.Hyperlinks.Add($ws.Range("D47"), "https://leetcode.com/problems/sum-of-nodes-with-even-valued-grandparent/")
$ws.Range("E47").Value = "recursion"
$ws.Range("F47").Value = "Medium"
$ws.Range("G47").Value = "No,See how I approached the problem"
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = "Simple recursion, see if current node value is even , add sum of all the grad child values else mode to the left and right child"

# ---------------------------------------------------------------------
# Row 48: "Count Square Submatrices with All Ones"
# ---------------------------------------------------------------------
$ws.Range("A48").Value = 44098
$ws.Range("B48").Value = 1277
$ws.Range("C48").Value = 1
$ws.Range("E48").Value = "dp"
$ws.Range("F48").Value = "Medium"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = @"
#no of possible square that end at point i,j will
#be equal to the min of max possible square out of
# square ending at i-1,j-1 & i-1,j & i,j-1
                    
"@
$ws.Range("G48").Value = "yes , 1 time"
$ws.Range("D48").Value = "https://leetcode.com/problems/count-square-submatrices-with-all-ones/submissions/"
$ws.Hyperlinks.Add($ws.Range("D48"), "https://leetcode.com/problems/count-square-submatrices-with-all-ones/submissions/")

# ---------------------------------------------------------------------
# Match formatting of the previous last row (46) for both new rows -
# reuses the existing date / body / hyperlink cell styles instead of
# minting new ones.
# ---------------------------------------------------------------------
$ws.Range("A46:J46").Copy()
$ws.Range("A47:J48").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row heights: 30pt for row 47 (single line), 60pt for row 48 (wrapped
# multi-line note).
$ws.Rows.Item(47).RowHeight = 30
$ws.Rows.Item(48).RowHeight = 60

# ---------------------------------------------------------------------
# View: scroll down a couple rows and select the freshly added last row
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("A48:J48").Select()

Write-Output "Added rows 47-48 to 'Questions Tracking Sheet'"
